$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 83 (the "「心からの言葉」" post), shifting all
# subsequent rows up by one. This matches the commit which removed
# that blog post entry from the sheet.
$ws.Rows.Item(83).Delete()
